$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing last row (row 135), pushing it to row 136.
# The new row inherits formatting (including the date style) from the row it
# was inserted at.
$ws.Rows("135").Insert()

# Populate the newly inserted row 135 with this week's new data.
$ws.Range("A135").Value = 3
$ws.Range("B135").Value = "Femacal de La Calera"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 44595
$ws.Range("E135").Value = 5
$ws.Range("F135").Value = 100112030
$ws.Range("G135").Value = "Poroto granado"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 73
$ws.Range("K135").Value = 22000
$ws.Range("L135").Value = 23000
$ws.Range("M135").Value = 22521
$ws.Range("N135").Value = '$/malla 25 kilos'
$ws.Range("O135").Value = "Provincia de Quillota"
$ws.Range("P135").Value = 901
$ws.Range("Q135").Value = 25
$ws.Range("R135").Value = "Hortaliza"
